$d = $word.ActiveDocument
$target = "Invalid block: Unexpected tag EOF missing [ENDFOR] while parsing m:for v| self.eClassifiers"
$prefix = "    <---"
$new = $prefix + $target
$targetLen = $target.Length

# Collect the start position of every occurrence first (single shared Range object).
$r = $d.Content
$find = $r.Find
$find.ClearFormatting()
$starts = New-Object System.Collections.ArrayList
$found = $find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    [void]$starts.Add($r.Start)
    $r.Collapse(0)
    [void]$r.MoveEnd(1, 1000000)
    $found = $find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# Enable track changes so replacing text re-uses / splits runs precisely at
# the edited span instead of Word re-coalescing same-format runs together.
$d.TrackRevisions = $true

# Replace from the last occurrence to the first so earlier offsets stay valid.
for ($i = $starts.Count - 1; $i -ge 0; $i--) {
    $pos = $starts[$i]
    $editRange = $d.Range($pos, $pos + $targetLen)
    $editRange.Text = $new
}

# Clean up: fold the edits back into normal (non-tracked) content.
$d.AcceptAllRevisions()
$d.TrackRevisions = $false

Write-Host "Updated $($starts.Count) occurrences"
